# Atualização de bases das ligas, do dia: 13-06-2024 às 19:35
#
# The underlying data rows (B:AD) for several fixtures were stored against
# the wrong row (the sequential index in column A and the date in column D
# stayed put, but the match id / teams / odds data had been shuffled).
# This script rotates the B:AD payload of each affected row-group back into
# the correct row, leaving column A (row index) and column D (date) intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is a group of worksheet rows whose B:AD content must be
# rotated by one position: row[i] takes the B:AD payload that currently
# lives on row[i+1] (wrapping around).
$groups = @(
    @(9, 10, 11),
    @(52, 54),
    @(126, 127),
    @(142, 143),
    @(147, 148),
    @(201, 202, 203, 204)
)

foreach ($group in $groups) {
    # Snapshot the current B:AD values for every row in the group before
    # writing anything back, since the rotation is circular.
    $snapshots = @{}
    foreach ($r in $group) {
        $snapshots[$r] = $ws.Range("B$r`:AD$r").Value2
    }

    $count = $group.Count
    for ($i = 0; $i -lt $count; $i++) {
        $targetRow = $group[$i]
        $sourceRow = $group[($i + 1) % $count]
        $ws.Range("B$targetRow`:AD$targetRow").Value2 = $snapshots[$sourceRow]
    }
}
